# Apply the "gh-pages output regenerated" update to 杭州-漫展信息.xlsx
# - refreshes "want to go" counters (column F) for a number of existing
#   rows on both the "展览" and "全部类型" sheets
# - on "全部类型" a newly scraped event (AD02动漫展--钟晨瑶内场票) is
#   prepended to the 2024-03/04 block, pushing rows 37-44 down to 38-45
#   and dropping the previous last row (赛马娘only) off the bottom

$wb = $excel.ActiveWorkbook

# column B holds plain "yyyy-mm-dd" text. Assigning a bare string like that
# through .Value gets auto-detected as a real date by Excel, so force the
# cell to Text first and strip the formatting mark afterwards to land back
# on a plain string cell (matches how the source file stores it).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# Sheet "展览": simple refreshed counts in column F
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

$expoUpdates = @{
    6  = 5143
    7  = 421
    15 = 16
    17 = 1700
    18 = 1437
    19 = 786
    22 = 289
    23 = 493
    25 = 1043
    28 = 2326
    31 = 70
    39 = 606
}
foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Range("F$row").Value = $expoUpdates[$row]
}

# ---------------------------------------------------------------------
# Sheet "全部类型": same refreshed counts in column F (rows 7-34)
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$allUpdates = @{
    7  = 5143
    8  = 421
    20 = 16
    23 = 1700
    24 = 1437
    25 = 786
    28 = 289
    30 = 493
    32 = 1043
    34 = 2326
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}

# ---------------------------------------------------------------------
# Sheet "全部类型": push rows 37-44 down into 38-45 (keep column A as-is)
# then write the brand-new row 37. This mirrors a fresh scrape being
# prepended to the list while the tail entry (赛马娘only, old row 45)
# falls off the bottom.
# ---------------------------------------------------------------------
for ($row = 44; $row -ge 37; $row--) {
    $srcRow = $row
    $dstRow = $row + 1
    Set-TextValue $wsAll.Range("B$dstRow") $wsAll.Range("B$srcRow").Value()
    $wsAll.Range("C$dstRow").Value = $wsAll.Range("C$srcRow").Value()
    $wsAll.Range("D$dstRow").Value = $wsAll.Range("D$srcRow").Value()
    $wsAll.Range("E$dstRow").Value = $wsAll.Range("E$srcRow").Value()
    $wsAll.Range("F$dstRow").Value = $wsAll.Range("F$srcRow").Value()
    $wsAll.Range("G$dstRow").Value = $wsAll.Range("G$srcRow").Value()
    $wsAll.Range("H$dstRow").Value = $wsAll.Range("H$srcRow").Value()
    $wsAll.Range("I$dstRow").Value = $wsAll.Range("I$srcRow").Value()
}

# the row that lands on 44 (originally row 43, "倒霉死勒内场票") carries an
# already-refreshed want-to-go count
$wsAll.Range("F44").Value = 606

# brand-new row 37
Set-TextValue $wsAll.Range("B37") "2024-03-24"
$wsAll.Range("C37").Value = "杭州·AD02动漫展--钟晨瑶内场票"
$wsAll.Range("D37").Value = "钱江世纪城奔竞大道353号 杭州国际博览中心"
$wsAll.Range("E37").Value = "2024.03.24 09:30-03.24 17:00"
$wsAll.Range("F37").Value = 70
$wsAll.Range("G37").Value = 258
$wsAll.Range("H37").Value = "https://show.bilibili.com/platform/detail.html?id=81820"
$wsAll.Range("I37").Value = "//i1.hdslb.com/bfs/openplatform/202402/aHRmCxr31707296105225.jpeg"
